$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194099426269531
$ws.Range("B1").Value = 1.212843298912048
$ws.Range("C1").Value = 6.823255062103271
$ws.Range("D1").Value = 2.099967956542969
$ws.Range("E1").Value = 1.147705078125
